# Apply the shelter-reassignment updates (renamed shelterData to shelData)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 6: Lias -> now assigned to Ibayo Elementary School
$ws.Range("D6").Value = "Ibayo Elementary School"
$ws.Range("E6").Value = 120.959816737558
$ws.Range("F6").Value = 14.7535649557989

# Row 8: Nagbalon -> now assigned to Old Municipal Bldg.
$ws.Range("D8").Value = "Old Municipal Bldg."
$ws.Range("E8").Value = 120.948177254006
$ws.Range("F8").Value = 14.7573006861396

# Row 9: Patubig -> now assigned to Marilao Central School
$ws.Range("D9").Value = "Marilao Central School"
$ws.Range("E9").Value = 120.949191076043
$ws.Range("F9").Value = 14.7549081782114

# Row 10: Poblacion I -> now assigned to Barangay Hall Nagbalon
$ws.Range("D10").Value = "Barangay Hall Nagbalon"
$ws.Range("E10").Value = 120.950788291388
$ws.Range("F10").Value = 14.7523618894178
